# "Generate Report for Archive"
#
# The status "Ready for handoff" moves on to "In Translation" for the file
# tracked in this report, and the Status-ish columns that held that text
# re-autosize (narrower, since "In Translation" is shorter) on the
# Overview sheet (columns E and F) as well as on the per-locale sheets
# (column C on "zh-cn" and "de-de").

$wb = $excel.ActiveWorkbook

# --- Update the status text everywhere it appears ---------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Re-size the columns that carried that text ------------------------
# ColumnWidth is expressed in "characters" and gets rounded to the
# nearest 1/6th of a character by the host; 12.5 is the closest
# achievable value that lands on the new target column width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
